$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.001.93"
$ws.Range("E2").Value = "  +10.63%  "
$ws.Range("D3").Value = "1.814.78"
$ws.Range("E3").Value = "  +7.44%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.541"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.79%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.97"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.96"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.279"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0669"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0928"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.32%  "
$ws.Range("D13").Value = "2.078.40"
$ws.Range("E13").Value = "  +7.53%  "
$ws.Range("D14").Value = "1.821.11"
$ws.Range("E14").Value = "  +7.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.640"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "34.012.64"
$ws.Range("E16").Value = "  +10.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "10.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "255.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.12%  "
$ws.Range("D21").Value = "0.0₃0743"
$ws.Range("E21").Value = "  +3.48%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.115"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.24%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.65%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.94%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0510"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.74%  "
$ws.Range("D35").Value = "1.556.38"
$ws.Range("E35").Value = "  +2.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.51%  "
$ws.Range("E37").Value = "  +4.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0187"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "84.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.617"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.905"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0523"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.61%  "
$ws.Range("E46").Value = "  +3.80%  "
$ws.Range("D47").Value = "1.968.68"
$ws.Range("E47").Value = "  +7.55%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  +7.87%  "
